$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 7, shifting rows 7+ down by one
$ws.Rows("7").Insert()

$ws.Range("B7").Value = "Address"
$ws.Range("C7").Value = "adr"

$ws.Range("C8").Select()
